# Actualizacion automatica de scrims_actualizado.xlsx (2025-07-25 18:02:44)
# Adds new scrim result rows to three worksheets, extending each sheet's
# used range (dimension) and copying the formatting of the most recent
# existing row so the new rows look consistent with the rest of the table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Crystal Arcade": append rows 52-54 (dimension A3:N51 -> A3:N54)
# ---------------------------------------------------------------------
$wsCA = $wb.Worksheets.Item("Crystal Arcade")

# Row 52 - same "Equipo 1" look as row 51
$wsCA.Range("A51:N51").Copy($wsCA.Range("A52:N52"))
$wsCA.Range("A52").Value = "JANET"
$wsCA.Range("B52").Value = "BUZZ"
$wsCA.Range("C52").Value = "CORDELIUS"
$wsCA.Range("D52").Value = "ASH"
$wsCA.Range("E52").Value = "KIT"
$wsCA.Range("F52").Value = "GUS"
$wsCA.Range("G52").Value = "Equipo 1"
$wsCA.Range("H52").Value = "FUT|GeRo"
$wsCA.Range("I52").Value = "FUT|Nowy297"
$wsCA.Range("J52").Value = "FUT|MeOw"
$wsCA.Range("K52").Value = "TH|LeNain"
$wsCA.Range("L52").Value = "TH|iKaoss"
$wsCA.Range("M52").Value = "TH|Zhar"
$wsCA.Range("N52").Value = "20250725T153816.000Z"

# Row 53 - "Equipo 2" look, borrow the G-cell style from another sheet
# that already has a row in the newer "Equipo 2" format.
$wsCA.Range("A51:N51").Copy($wsCA.Range("A53:N53"))
$wsDS = $wb.Worksheets.Item("Dry Season")
$wsDS.Range("G42").Copy($wsCA.Range("G53"))
$wsCA.Range("A53").Value = "JANET"
$wsCA.Range("B53").Value = "BUZZ"
$wsCA.Range("C53").Value = "CORDELIUS"
$wsCA.Range("D53").Value = "ASH"
$wsCA.Range("E53").Value = "KIT"
$wsCA.Range("F53").Value = "GUS"
$wsCA.Range("G53").Value = "Equipo 2"
$wsCA.Range("H53").Value = "FUT|GeRo"
$wsCA.Range("I53").Value = "FUT|Nowy297"
$wsCA.Range("J53").Value = "FUT|MeOw"
$wsCA.Range("K53").Value = "TH|LeNain"
$wsCA.Range("L53").Value = "TH|iKaoss"
$wsCA.Range("M53").Value = "TH|Zhar"
$wsCA.Range("N53").Value = "20250725T153601.000Z"

# Row 54 - same "Equipo 1" look as row 51
$wsCA.Range("A51:N51").Copy($wsCA.Range("A54:N54"))
$wsCA.Range("A54").Value = "JANET"
$wsCA.Range("B54").Value = "BUZZ"
$wsCA.Range("C54").Value = "CORDELIUS"
$wsCA.Range("D54").Value = "ASH"
$wsCA.Range("E54").Value = "KIT"
$wsCA.Range("F54").Value = "GUS"
$wsCA.Range("G54").Value = "Equipo 1"
$wsCA.Range("H54").Value = "FUT|GeRo"
$wsCA.Range("I54").Value = "FUT|Nowy297"
$wsCA.Range("J54").Value = "FUT|MeOw"
$wsCA.Range("K54").Value = "TH|LeNain"
$wsCA.Range("L54").Value = "TH|iKaoss"
$wsCA.Range("M54").Value = "TH|Zhar"
$wsCA.Range("N54").Value = "20250725T153305.000Z"

# ---------------------------------------------------------------------
# Sheet "Dry Season": append rows 43-44 (dimension A3:N42 -> A3:N44)
# ---------------------------------------------------------------------
$wsDS.Range("A42:N42").Copy($wsDS.Range("A43:N43"))
$wsDS.Range("A43").Value = "BROCK"
$wsDS.Range("B43").Value = "FANG"
$wsDS.Range("C43").Value = "GENE"
$wsDS.Range("D43").Value = "CARL"
$wsDS.Range("E43").Value = "SQUEAK"
$wsDS.Range("F43").Value = "BELLE"
$wsDS.Range("G43").Value = "Equipo 2"
$wsDS.Range("H43").Value = "BBO|Adrii"
$wsDS.Range("I43").Value = "BBO|Bogdan"
$wsDS.Range("J43").Value = "SUP|Salty"
$wsDS.Range("K43").Value = "NXT|Rup"
$wsDS.Range("L43").Value = "NXT|Arthur"
$wsDS.Range("M43").Value = "NXT|amos"
$wsDS.Range("N43").Value = "20250725T153419.000Z"

$wsDS.Range("A42:N42").Copy($wsDS.Range("A44:N44"))
$wsDS.Range("A44").Value = "BROCK"
$wsDS.Range("B44").Value = "FANG"
$wsDS.Range("C44").Value = "GENE"
$wsDS.Range("D44").Value = "CARL"
$wsDS.Range("E44").Value = "SQUEAK"
$wsDS.Range("F44").Value = "BELLE"
$wsDS.Range("G44").Value = "Equipo 2"
$wsDS.Range("H44").Value = "BBO|Adrii"
$wsDS.Range("I44").Value = "BBO|Bogdan"
$wsDS.Range("J44").Value = "SUP|Salty"
$wsDS.Range("K44").Value = "NXT|Rup"
$wsDS.Range("L44").Value = "NXT|Arthur"
$wsDS.Range("M44").Value = "NXT|amos"
$wsDS.Range("N44").Value = "20250725T153158.000Z"

# ---------------------------------------------------------------------
# Sheet "Pit Stop": append rows 55-56 (dimension A3:N54 -> A3:N56)
# ---------------------------------------------------------------------
$wsSF = $wb.Worksheets.Item("Pit Stop")

$wsSF.Range("A53:N53").Copy($wsSF.Range("A55:N55"))
$wsSF.Range("A55").Value = "HANK"
$wsSF.Range("B55").Value = "SHADE"
$wsSF.Range("C55").Value = "LUMI"
$wsSF.Range("D55").Value = "NITA"
$wsSF.Range("E55").Value = "BULL"
$wsSF.Range("F55").Value = "MICO"
$wsSF.Range("G55").Value = "Equipo 2"
$wsSF.Range("H55").Value = "NOVO|Subeme"
$wsSF.Range("I55").Value = "NOVO|Marco"
$wsSF.Range("J55").Value = "NOVO|Biso"
$wsSF.Range("K55").Value = "IC|Mebius"
$wsSF.Range("L55").Value = "IC|Nob?"
$wsSF.Range("M55").Value = "IC|RamaZR"
$wsSF.Range("N55").Value = "20250725T153226.000Z"

$wsSF.Range("A53:N53").Copy($wsSF.Range("A56:N56"))
$wsSF.Range("A56").Value = "HANK"
$wsSF.Range("B56").Value = "SHADE"
$wsSF.Range("C56").Value = "LUMI"
$wsSF.Range("D56").Value = "NITA"
$wsSF.Range("E56").Value = "BULL"
$wsSF.Range("F56").Value = "MICO"
$wsSF.Range("G56").Value = "Equipo 2"
$wsSF.Range("H56").Value = "NOVO|Subeme"
$wsSF.Range("I56").Value = "NOVO|Marco"
$wsSF.Range("J56").Value = "NOVO|Biso"
$wsSF.Range("K56").Value = "IC|Mebius"
$wsSF.Range("L56").Value = "IC|Nob?"
$wsSF.Range("M56").Value = "IC|RamaZR"
$wsSF.Range("N56").Value = "20250725T153043.000Z"
